$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "01.09.21.50UM.003_PROCESSED.xlsx"
$ws.Range("D2").Value = 68.5745443015133
$ws.Range("E2").Value = 50.44337779213091
$ws.Range("F2").Value = 0.4824206896627805
$ws.Range("G2").Value = 0.0007174634465868938
$ws.Range("H2").Value = -247353.7475038553
$ws.Range("I2").Value = 2280.290147728559
$ws.Range("J2").Value = 788561.0695067085
$ws.Range("K2").Value = 7.43478449457637
$ws.Range("L2").Value = 6.245959085320488
$ws.Range("M2").Value = 7.990934214151545
$ws.Range("N2").Value = 0.000001
$ws.Range("R2").Value = 0.001853821505849047
